# Saldo.xlsx update: add/remove rows in the "Export" sheet.
# Process bottom-to-top so the row numbers used below stay valid for
# steps still to come.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force the value to be stored as text (preserves leading zeros in
    # account numbers, like the existing "Conta" column cells), then
    # drop the temporary "@" number format so the cell keeps the same
    # (default) style as its neighbours.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# 1) Remove the row for 005685089 / CARNEIRO / 1704.16 (currently row 16)
$ws.Rows.Item(16).Delete()

# 2) Insert a new row before 004368468 / AHMAD (currently row 14) for
#    005073033 / NILBORN / 3983.35
$ws.Rows.Item(14).Insert()
Set-TextCell 14 1 "005073033"
Set-TextCell 14 2 "NILBORN"
$ws.Cells.Item(14, 3).Value = 3983.35

# 3) Remove the row for 008243633 / DANIELA / 8300 (currently row 12)
$ws.Rows.Item(12).Delete()

# 4) Insert two new rows before 005002390 / LUCIANO (currently row 11):
#    004472431 / LUIS / 9952.64
#    004404248 / PAULO / 9951.05
$ws.Rows.Item(11).Insert()
Set-TextCell 11 1 "004472431"
Set-TextCell 11 2 "LUIS"
$ws.Cells.Item(11, 3).Value = 9952.64

$ws.Rows.Item(12).Insert()
Set-TextCell 12 1 "004404248"
Set-TextCell 12 2 "PAULO"
$ws.Cells.Item(12, 3).Value = 9951.05

# 5) Insert a new row before 001761119 / BLUEMETRIX (currently row 4):
#    004474776 / GILSON / 49764.44
$ws.Rows.Item(4).Insert()
Set-TextCell 4 1 "004474776"
Set-TextCell 4 2 "GILSON"
$ws.Cells.Item(4, 3).Value = 49764.44
